$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.765.74'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '1.540.27'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3881'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3190'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.79%  '
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07208'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('E11').Value = '  -6.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.648'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.61'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.616'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001116'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.544.68'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06589'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.158'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.384'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.18%  '
$ws.Range('D25').Value = '21.772.64'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.401'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '146.59'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.841'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('D30').Value = '1.719.88'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.64'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9742'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -12.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.915'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08202'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.948'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.32%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.167'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06087'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.483'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -17.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02211'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2041'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.196'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.72'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5775'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.44%  '
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.748'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5535'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '118.37'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.873'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.147'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06728'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.59%  '
